$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bump the suite version from 0.1 to 1.0 (keep it stored as text, like the original "0.1")
$ws.Range("D2").Value = "'1.0"

# 2. TC2's second step (cancel a diária) and TC3's second step (filter by params)
#    had their contents swapped - this is effectively a new TC3 ("filter search")
#    being inserted ahead of the "cancel" scenario, which moved to TC3 in its place.
$ws.Range("B20").Value = "Chefe Indica alguns parâmetros específicos para a busca; Informa o nome do beneficiário; Filtra a listagem de solicitações."
$ws.Range("D20").Value = "SYSTEM Exibe uma nova listagem de solicitações, de acordo com os filtros informados pelo usuário."

$ws.Range("B28").Value = "Chefe Clica para realizar o cancelamento de uma diária."
$ws.Range("D28").Value = "SYSTEM Verifica que a solicitação está em situação SOLICITADA; Exibe mensagem de confirmação (MSG987 - Cancelar solicitação de diária) para o usuário (que deve confirmar); Cancela a diária, mudando sua situação para CANCELADA (ver diagrama de estados da diária)."
